$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 457, pushing old rows 457:538 down to 458:539.
$ws.Rows("457:457").Insert()

# Populate the newly inserted row 457. It mirrors the row that is now at 458
# (the former row 457) except for the "Fecha" (D) and "Volumen" (J) values.
$ws.Range("A457").Value = 10
$ws.Range("B457").Value = "Vega Modelo de Temuco"
$ws.Range("C457").Value = "La Araucanía"
$ws.Range("D457").Value = 45015
$ws.Range("E457").Value = 9
$ws.Range("F457").Value = 100112040
$ws.Range("G457").Value = "Cilantro"
$ws.Range("H457").Value = "Sin especificar"
$ws.Range("I457").Value = "Primera"
$ws.Range("J457").Value = 100
$ws.Range("K457").Value = 5000
$ws.Range("L457").Value = 5000
$ws.Range("M457").Value = 5000
$ws.Range("N457").Value = "$/docena de atados (2 kilos)"
$ws.Range("O457").Value = "Provincia de Cautín"
$ws.Range("P457").Value = 2500
$ws.Range("Q457").Value = 2
$ws.Range("R457").Value = "Hortaliza"

# Re-apply the date number format (style index 2 in the original workbook)
# to the new D457 cell, matching the rest of the "Fecha" column.
$ws.Range("D457").NumberFormat = "YYYY-MM-DD HH:MM:SS"
